$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74 - PEL.NS
$ws.Range("A74").Value = "PEL.NS"
$ws.Range("B74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B74").Value = 37803
$ws.Range("C74").Value = 30
$ws.Range("D74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D74").Value = 37073
$ws.Range("E74").Value = 22
$ws.Range("F74").Value = 21.93608474731445
$ws.Range("G74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G74").Value = 37438
$ws.Range("H74").Value = 26
$ws.Range("I74").Value = 21.79438591003418
$ws.Range("J74").Value = "Low"
$ws.Range("K74").Value = -0.03542470932006836
$ws.Range("L74").Value = 22.71542835235596
$ws.Range("M74").Value = 3
$ws.Range("N74").Value = 1
$ws.Range("O74").Value = 2

# Row 75 - PEL.NS
$ws.Range("A75").Value = "PEL.NS"
$ws.Range("B75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B75").Value = 44835
$ws.Range("C75").Value = 107
$ws.Range("D75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D75").Value = 42826
$ws.Range("E75").Value = 85
$ws.Range("F75").Value = 1741.661987304688
$ws.Range("G75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G75").Value = 44470
$ws.Range("H75").Value = 103
$ws.Range("I75").Value = 1751.69140625
$ws.Range("J75").Value = "High"
$ws.Range("K75").Value = 0.55718994140625
$ws.Range("L75").Value = 1694.300842285156
$ws.Range("M75").Value = 3
$ws.Range("N75").Value = 1
$ws.Range("O75").Value = 2
